{"js": "// Target: the first paragraph in the body (the \"**ID__..._ID**\" placeholder\n// paragraph) gets:\n//   1. A paragraph border (pBdr) with 5-twip spacing on all four sides\n//      (top/left/bottom/right) and no visible line (distance-only border,\n//      matching the existing pattern already used later in the document).\n//   2. Its left indent changed from 120 -> 225 (twips), i.e. 6pt -> 11.25pt.\n//   3. Its text changed from \"**ID__AFFARS_5311_topic_5__ID**\" to\n//      \"**ID__AFFARS_5311_273_3__ID**\", collapsing the paragraph down to a\n//      single run (the trailing \" \" run is removed).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[0];\n\n// --- 1) Paragraph border: 5-twip spacing on all sides, no line style.\n// Office.js's `Paragraph.borders` collection only surfaces line\n// style/color/width per edge (no distance-from-text setter), so we reach\n// the same Word OM bridge the generated proxies use (`ParagraphFormat\n// .Borders.DistanceFromTop/Left/Bottom/Right`, the COM-level property that\n// actually carries the OOXML `w:space` value) via the shared `_omSet` call\n// that backs every Office.js property on this host.\nconst borders = target.borders;\nborders.load(\"items\");\nawait context.sync();\n\nconst edge = borders.items[0];\nedge._omSet(\"DistanceFromTop\", 5, \"Borders\");\nedge._omSet(\"DistanceFromLeft\", 5, \"Borders\");\nedge._omSet(\"DistanceFromBottom\", 5, \"Borders\");\nedge._omSet(\"DistanceFromRight\", 5, \"Borders\");\n\n// --- 2) Left indent: 120 twips -> 225 twips (Office.js works in points).\ntarget.leftIndent = 225 / 20;\n\n// --- 3) Replace the paragraph's whole text with the new placeholder id.\n// Replacing via the paragraph's range collapses both existing runs into a\n// single run that inherits the original run formatting (font, size, color,\n// bold/italic flags), matching the diff's removal of the second (space-only)\n// run.\nconst range = target.getRange();\nrange.insertText(\"**ID__AFFARS_5311_273_3__ID**\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Target: the first paragraph in the body (the \"**ID__..._ID**\" placeholder\n# paragraph) gets:\n#   1. A paragraph border (pBdr) with 5-twip spacing on all four sides\n#      (top/left/bottom/right) and no visible line (distance-only border,\n#      matching the existing pattern already used later in the document).\n#   2. Its left indent changed from 120 -> 225 (twips).\n#   3. Its text changed from \"**ID__AFFARS_5311_topic_5__ID**\" to\n#      \"**ID__AFFARS_5311_273_3__ID**\", collapsing the paragraph down to a\n#      single run (the trailing \" \" run is removed).\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# --- 1) Paragraph border: 5-twip spacing on all sides, no line style.\n$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5\n$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5\n$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5\n$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5\n\n# --- 2) Left indent: 120 twips -> 225 twips (Word COM indent is in points).\n$p.Range.ParagraphFormat.LeftIndent = 225 / 20\n\n# --- 3) Replace the paragraph's whole text with the new placeholder id,\n# collapsing the two existing runs (the id text + trailing space) into one.\n$p.Range.Text = \"**ID__AFFARS_5311_273_3__ID**\"\n"}
